$wb = $excel.ActiveWorkbook

# --- Sheet 'Recommandations' ---
$ws1 = $wb.Worksheets.Item("Recommandations")

$data1 = @(
    ,@("UNIWAX CI", 0, 4, 3425, 915, "🟡 Observer", "➖ Neutre")
    ,@("BRVM - SERVICES PUBLICS", 0, 8, 3308.61, 103.71, "🟡 Observer", "➖ Neutre")
    ,@("CFAO MOTORS CI", 0, 4, 3245, 910, "🟡 Observer", "➖ Neutre")
    ,@("SETAO CI", 0, 4, 2750, 725, "🟡 Observer", "➖ Neutre")
    ,@("BRVM - AUTRES SECTEURS", 0, 4, 2524.52, 629.42, "🟡 Observer", "➖ Neutre")
    ,@("NEI-CEDA CI", 0, 4, 2505, 645, "🟡 Observer", "➖ Neutre")
    ,@("AIR LIQUIDE CI", 0, 4, 2380, 600, "🟡 Observer", "➖ Neutre")
    ,@("BRVM - DISTRIBUTION", 0, 4, 1559.41, 398.72, "🟡 Observer", "➖ Neutre")
    ,@("BRVM - TRANSPORT", 0, 4, 1386.63, 348.8, "🟡 Observer", "➖ Neutre")
    ,@("BRVM - AGRICULTURE", 0, 4, 1320.82, 329.57, "🟡 Observer", "➖ Neutre")
    ,@("BRVM - INDUSTRIELS", 0, 4, 552.3, 137.76, "🟡 Observer", "➖ Neutre")
    ,@("BRVM-PRESTIGE", 0, 4, 530.95, 132.79, "🟡 Observer", "➖ Neutre")
    ,@("BRVM - FINANCES", 0, 4, 501.07, 125.97, "🟡 Observer", "➖ Neutre")
    ,@("BRVM - SERVICES FINANCIERS", 0, 4, 492.44, 123.8, "🟡 Observer", "➖ Neutre")
    ,@("BRVM - CONSOMMATION DISCRETIONNAIRE", 0, 4, 474.71, 124.73, "🟡 Observer", "➖ Neutre")
    ,@("BRVM - ENERGIE", 0, 4, 438.18, 108.65, "🟡 Observer", "➖ Neutre")
    ,@("BRVM - TELECOMMUNICATIONS", 0, 4, 377.91, 94.34, "🟡 Observer", "➖ Neutre")
    ,@("BRVM - INDUSTRIE       (**)", 0, 1, 262.79, 262.79, "🟡 Observer", "➖ Neutre")
    ,@("BRVM - INDUSTRIE                  (**)", 0, 1, 251.95, 251.95, "🟡 Observer", "➖ Neutre")
    ,@("BRVM - CONSOMMATION DE BASE      (**)", 0, 1, 217.36, 217.36, "🟡 Observer", "➖ Neutre")
    ,@("BRVM - CONSOMMATION DE BASE          (**)", 0, 1, 209.51, 209.51, "🟡 Observer", "➖ Neutre")
    ,@("BRVM-PRINCIPAL          (**)", 0, 1, 192.79, 192.79, "🟡 Observer", "➖ Neutre")
    ,@("BRVM-PRINCIPAL                    (**)", 0, 1, 191.6, 191.6, "🟡 Observer", "➖ Neutre")
    ,@("CFAO MOTORS CI (CFAC)", 4, 0, 28.66, 7.06, "🟢 Achat", "✅ Renforcer")
    ,@("UNIWAX CI (UNXC)", 3, 0, 19.17, 4.57, "🟢 Achat", "✅ Renforcer")
    ,@("NEI-CEDA CI (NEIC)", 2, 0, 11.15, 6.98, "🟡 Observer", "➖ Neutre")
    ,@("ORAGROUP TOGO (ORGT)", 2, 1, 5.39, 5.57, "🟡 Observer", "👀 À surveiller")
    ,@("SETAO CI (STAC)", 1, 1, 4.61, -2.74, "🟡 Observer", "👀 À surveiller")
    ,@("TOTALENERGIES MARKETING CI (TTLC)", 1, 0, 4.26, 4.26, "🟡 Observer", "➖ Neutre")
    ,@("AIR LIQUIDE CI (SIVC)", 1, 0, 3.48, 3.48, "🟡 Observer", "➖ Neutre")
    ,@("BANK OF AFRICA ML (BOAM)", 1, 0, 2.23, 2.23, "🟡 Observer", "➖ Neutre")
    ,@("SOCIETE GENERALE COTE D'IVOIRE (SGBC)", 1, 0, 1.96, 1.96, "🟡 Observer", "➖ Neutre")
    ,@("SOGB CI (SOGC)", 1, 1, 0.39, -2.47, "🟡 Observer", "👀 À surveiller")
    ,@("ECOBANK TRANS. INCORP. TG (ETIT)", 0, 1, 0.37, 6.25, "🟡 Observer", "👀 À surveiller")
    ,@("TOTAL", 0, 3, 0, 0, "🟡 Observer", "➖ Neutre")
    ,@("TOTALENERGIES MARKETING SN (TTLS)", 1, 1, -0.27, 3.2, "🟡 Observer", "👀 À surveiller")
    ,@("BANK OF AFRICA BF (BOABF)", 0, 1, -1.05, -1.05, "🟡 Observer", "➖ Neutre")
    ,@("AFRICA GLOBAL LOGISTICS CI (SDSC)", 0, 1, -1.06, -1.06, "🟡 Observer", "➖ Neutre")
    ,@("SICABLE CI (CABC)", 0, 1, -1.13, -1.13, "🟡 Observer", "➖ Neutre")
    ,@("FILTISAC CI (FTSC)", 0, 1, -1.24, -1.24, "🟡 Observer", "➖ Neutre")
    ,@("VIVO ENERGY CI (SHEC)", 0, 1, -1.3, -1.3, "🟡 Observer", "➖ Neutre")
    ,@("ONATEL BF (ONTBF)", 0, 1, -2.51, -2.51, "🟡 Observer", "➖ Neutre")
    ,@("SITAB CI (STBC)", 0, 1, -3.04, -3.04, "🟡 Observer", "➖ Neutre")
    ,@("SUCRIVOIRE (SCRC)", 0, 1, -3.28, -3.28, "🟡 Observer", "➖ Neutre")
    ,@("SERVAIR ABIDJAN CI (ABJC)", 0, 1, -3.31, -3.31, "🟡 Observer", "➖ Neutre")
    ,@("BANK OF AFRICA NG (BOAN)", 0, 1, -4.94, -4.94, "🟡 Observer", "➖ Neutre")
    ,@("BERNABE CI (BNBC)", 0, 3, -7.83, -1.96, "🔴 Vente", "⚠️ Risque de décrochage")
    ,@("UNILEVER CI (UNLC)", 0, 2, -14.46, -7.49, "🟡 Observer", "➖ Neutre")
)

for ($i = 0; $i -lt $data1.Length; $i++) {
    $row = $data1[$i]
    $r = $i + 2
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $ws1.Cells.Item($r, 6).Value = $row[5]
    $ws1.Cells.Item($r, 7).Value = $row[6]
}

# --- Sheet 'Top_YTD' ---
$ws2 = $wb.Worksheets.Item("Top_YTD")

$data2 = @(
    ,7857305
    ,828664.5
    ,678314.48
    ,383876.45
    ,285643.47
    ,277794.31
    ,233176.75
    ,57457.23
    ,39698.91
    ,34153.08
)

for ($i = 0; $i -lt $data2.Length; $i++) {
    $r = $i + 2
    $ws2.Cells.Item($r, 2).Value = $data2[$i]
}

